# Update 想去人数 (F column) values across all sheets to the new crawl totals.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1713
$ws.Cells.Item(4, 6).Value = 9947
$ws.Cells.Item(5, 6).Value = 9
$ws.Cells.Item(6, 6).Value = 294
$ws.Cells.Item(8, 6).Value = 91
$ws.Cells.Item(12, 6).Value = 69
$ws.Cells.Item(13, 6).Value = 1554
$ws.Cells.Item(15, 6).Value = 343
$ws.Cells.Item(18, 6).Value = 455
$ws.Cells.Item(19, 6).Value = 1144
$ws.Cells.Item(24, 6).Value = 323
$ws.Cells.Item(29, 6).Value = 669
$ws.Cells.Item(31, 6).Value = 26
$ws.Cells.Item(32, 6).Value = 204
$ws.Cells.Item(34, 6).Value = 226
$ws.Cells.Item(35, 6).Value = 206
$ws.Cells.Item(36, 6).Value = 356
$ws.Cells.Item(38, 6).Value = 475
$ws.Cells.Item(39, 6).Value = 692
$ws.Cells.Item(40, 6).Value = 520
$ws.Cells.Item(42, 6).Value = 787
$ws.Cells.Item(43, 6).Value = 357
$ws.Cells.Item(44, 6).Value = 314
$ws.Cells.Item(45, 6).Value = 345

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(18, 6).Value = 1054
$ws.Cells.Item(20, 6).Value = 455
$ws.Cells.Item(21, 6).Value = 1089
$ws.Cells.Item(22, 6).Value = 313
$ws.Cells.Item(31, 6).Value = 195
$ws.Cells.Item(34, 6).Value = 149
$ws.Cells.Item(41, 6).Value = 57

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 79
$ws.Cells.Item(4, 6).Value = 801
$ws.Cells.Item(6, 6).Value = 2474
$ws.Cells.Item(7, 6).Value = 3950
$ws.Cells.Item(8, 6).Value = 45
$ws.Cells.Item(10, 6).Value = 221
$ws.Cells.Item(11, 6).Value = 160

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1713
$ws.Cells.Item(4, 6).Value = 801
$ws.Cells.Item(5, 6).Value = 9947
$ws.Cells.Item(7, 6).Value = 9
$ws.Cells.Item(8, 6).Value = 3950
$ws.Cells.Item(9, 6).Value = 45
$ws.Cells.Item(10, 6).Value = 221
$ws.Cells.Item(11, 6).Value = 221
$ws.Cells.Item(12, 6).Value = 91
$ws.Cells.Item(15, 6).Value = 343
$ws.Cells.Item(18, 6).Value = 456
$ws.Cells.Item(19, 6).Value = 1144
$ws.Cells.Item(24, 6).Value = 1054
$ws.Cells.Item(25, 6).Value = 323
$ws.Cells.Item(29, 6).Value = 1089
$ws.Cells.Item(30, 6).Value = 669
$ws.Cells.Item(33, 6).Value = 26
$ws.Cells.Item(34, 6).Value = 205
$ws.Cells.Item(36, 6).Value = 206
$ws.Cells.Item(38, 6).Value = 356
$ws.Cells.Item(40, 6).Value = 475
$ws.Cells.Item(41, 6).Value = 195
$ws.Cells.Item(42, 6).Value = 692
$ws.Cells.Item(43, 6).Value = 520
$ws.Cells.Item(44, 6).Value = 788
$ws.Cells.Item(45, 6).Value = 357
$ws.Cells.Item(47, 6).Value = 345
$ws.Cells.Item(50, 6).Value = 57
